$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.1303
$ws.Range("B4").Value = 4.627600000000003
$ws.Range("D4").Value = -6.985700000000001

$ws.Range("B5").Value = 5.217199999999999

$ws.Range("A7").Value = -21.67460000000001

$ws.Range("B8").Value = 4.922599999999997

$ws.Range("D9").Value = -8.540800000000008

$ws.Range("A16").Value = -21.54000000000001
$ws.Range("B16").Value = 4.825800000000001

$ws.Range("D18").Value = -8.421599999999991
